$wb = $excel.ActiveWorkbook

# --- Resumen sheet: update Maximo (C2) ---
$wsResumen = $wb.Worksheets("Resumen")
$wsResumen.Range("C2").Value = 690.8448751216349

# --- Solucion sheet: randomized reassignment of Pedido / Salida columns ---
$wsSolucion = $wb.Worksheets("Solucion")
$solucionData = @(
    @("Pedido_38", "S001"),
    @("Pedido_10", "S025"),
    @("Pedido_43", "S041"),
    @("Pedido_31", "S065"),
    @("Pedido_60", "S005"),
    @("Pedido_1", "S069"),
    @("Pedido_63", "S045"),
    @("Pedido_12", "S002"),
    @("Pedido_45", "S029"),
    @("Pedido_79", "S042"),
    @("Pedido_76", "S066"),
    @("Pedido_33", "S046"),
    @("Pedido_2", "S026"),
    @("Pedido_36", "S006"),
    @("Pedido_23", "S070"),
    @("Pedido_50", "S043"),
    @("Pedido_26", "S003"),
    @("Pedido_80", "S030"),
    @("Pedido_7", "S067"),
    @("Pedido_25", "S047"),
    @("Pedido_15", "S027"),
    @("Pedido_77", "S071"),
    @("Pedido_9", "S007"),
    @("Pedido_70", "S044"),
    @("Pedido_48", "S031"),
    @("Pedido_18", "S068"),
    @("Pedido_22", "S048"),
    @("Pedido_78", "S004"),
    @("Pedido_59", "S028"),
    @("Pedido_21", "S072"),
    @("Pedido_64", "S008"),
    @("Pedido_4", "S049"),
    @("Pedido_28", "S032"),
    @("Pedido_6", "S073"),
    @("Pedido_14", "S009"),
    @("Pedido_27", "S053"),
    @("Pedido_11", "S033"),
    @("Pedido_71", "S077"),
    @("Pedido_58", "S050"),
    @("Pedido_37", "S013"),
    @("Pedido_61", "S054"),
    @("Pedido_69", "S037"),
    @("Pedido_24", "S074"),
    @("Pedido_40", "S010"),
    @("Pedido_44", "S051"),
    @("Pedido_19", "S034"),
    @("Pedido_3", "S078"),
    @("Pedido_49", "S014"),
    @("Pedido_52", "S055"),
    @("Pedido_41", "S038"),
    @("Pedido_30", "S011"),
    @("Pedido_34", "S035"),
    @("Pedido_66", "S052"),
    @("Pedido_20", "S075"),
    @("Pedido_54", "S015"),
    @("Pedido_67", "S056"),
    @("Pedido_73", "S012"),
    @("Pedido_42", "S039"),
    @("Pedido_74", "S079"),
    @("Pedido_8", "S057"),
    @("Pedido_32", "S016"),
    @("Pedido_17", "S036"),
    @("Pedido_56", "S076"),
    @("Pedido_46", "S061"),
    @("Pedido_39", "S040"),
    @("Pedido_29", "S017"),
    @("Pedido_62", "S080"),
    @("Pedido_13", "S058"),
    @("Pedido_75", "S021"),
    @("Pedido_16", "S018"),
    @("Pedido_53", "S062"),
    @("Pedido_57", "S022"),
    @("Pedido_65", "S019"),
    @("Pedido_35", "S059"),
    @("Pedido_5", "S023"),
    @("Pedido_55", "S063"),
    @("Pedido_72", "S020"),
    @("Pedido_47", "S060"),
    @("Pedido_51", "S024"),
    @("Pedido_68", "S064")
)
for ($i = 0; $i -lt $solucionData.Count; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $solucionData[$i][0]
    $wsSolucion.Cells.Item($row, 2).Value = $solucionData[$i][1]
}

# --- Metricas sheet: update Tiempo per Zona ---
$wsMetricas = $wb.Worksheets("Metricas")
$wsMetricas.Range("B2").Value = 685.1088225754137
$wsMetricas.Range("B3").Value = 469.6813979889718
$wsMetricas.Range("B4").Value = 690.8448751216349
$wsMetricas.Range("B5").Value = 485.5336252567845
